{"js": "// Replace the date and each multiplication problem's text, one-to-one,\n// by searching for the exact old text and replacing it with the new text.\nconst replacements = [\n  [\"2024-01-30 Tuesday\", \"2024-01-31 Wednesday\"],\n  [\"922\u00d75=\", \"402\u00d74=\"],\n  [\"470\u00d72=\", \"347\u00d72=\"],\n  [\"104\u00d73=\", \"800\u00d76=\"],\n  [\"724\u00d79=\", \"295\u00d74=\"],\n  [\"761\u00d77=\", \"661\u00d74=\"],\n  [\"615\u00d78=\", \"308\u00d72=\"],\n  [\"380\u00d72=\", \"921\u00d75=\"],\n  [\"568\u00d73=\", \"923\u00d79=\"],\n  [\"176\u00d75=\", \"543\u00d75=\"],\n  [\"339\u00d75=\", \"680\u00d75=\"],\n  [\"258\u00d72=\", \"914\u00d77=\"],\n  [\"726\u00d77=\", \"242\u00d78=\"],\n  [\"205\u00d77=\", \"512\u00d72=\"],\n  [\"589\u00d75=\", \"902\u00d72=\"],\n  [\"959\u00d73=\", \"618\u00d74=\"],\n  [\"854\u00d75=\", \"355\u00d72=\"],\n  [\"189\u00d78=\", \"787\u00d74=\"],\n  [\"645\u00d77=\", \"438\u00d79=\"],\n  [\"786\u00d77=\", \"291\u00d73=\"],\n  [\"395\u00d75=\", \"624\u00d76=\"],\n  [\"426\u00d74=\", \"714\u00d75=\"],\n  [\"419\u00d74=\", \"163\u00d78=\"],\n  [\"484\u00d78=\", \"613\u00d73=\"],\n  [\"372\u00d74=\", \"692\u00d72=\"],\n  [\"487\u00d74=\", \"445\u00d78=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the date and each multiplication problem's text, one-to-one,\n# using Word's Find/Replace (Content.Find.Execute) against the whole document.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2024-01-30 Tuesday\", \"2024-01-31 Wednesday\")\n    ,@(\"922\u00d75=\", \"402\u00d74=\")\n    ,@(\"470\u00d72=\", \"347\u00d72=\")\n    ,@(\"104\u00d73=\", \"800\u00d76=\")\n    ,@(\"724\u00d79=\", \"295\u00d74=\")\n    ,@(\"761\u00d77=\", \"661\u00d74=\")\n    ,@(\"615\u00d78=\", \"308\u00d72=\")\n    ,@(\"380\u00d72=\", \"921\u00d75=\")\n    ,@(\"568\u00d73=\", \"923\u00d79=\")\n    ,@(\"176\u00d75=\", \"543\u00d75=\")\n    ,@(\"339\u00d75=\", \"680\u00d75=\")\n    ,@(\"258\u00d72=\", \"914\u00d77=\")\n    ,@(\"726\u00d77=\", \"242\u00d78=\")\n    ,@(\"205\u00d77=\", \"512\u00d72=\")\n    ,@(\"589\u00d75=\", \"902\u00d72=\")\n    ,@(\"959\u00d73=\", \"618\u00d74=\")\n    ,@(\"854\u00d75=\", \"355\u00d72=\")\n    ,@(\"189\u00d78=\", \"787\u00d74=\")\n    ,@(\"645\u00d77=\", \"438\u00d79=\")\n    ,@(\"786\u00d77=\", \"291\u00d73=\")\n    ,@(\"395\u00d75=\", \"624\u00d76=\")\n    ,@(\"426\u00d74=\", \"714\u00d75=\")\n    ,@(\"419\u00d74=\", \"163\u00d78=\")\n    ,@(\"484\u00d78=\", \"613\u00d73=\")\n    ,@(\"372\u00d74=\", \"692\u00d72=\")\n    ,@(\"487\u00d74=\", \"445\u00d78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
